# Apply updated shift-plan figures (rows 5-10) and update the saved selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 3
$ws.Range("P5").Value = 7
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 6
$ws.Range("X5").Value = 3

# Row 6
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 3
$ws.Range("P6").Value = 7
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = 7
$ws.Range("S6").Value = 6
$ws.Range("V6").Value = 6
$ws.Range("X6").Value = 3

# Row 7
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 3
$ws.Range("P7").Value = 7
$ws.Range("Q7").Value = 7
$ws.Range("R7").Value = 7
$ws.Range("S7").Value = 6
$ws.Range("X7").Value = 3

# Row 8
$ws.Range("S8").Value = 7
$ws.Range("T8").Value = 7

# Row 9
$ws.Range("O9").Value = 5
$ws.Range("X9").Value = 4

# Row 10
$ws.Range("M10").Value = 4
$ws.Range("V10").Value = 5

# Update the selected cell reflected in the saved sheet view
$ws.Range("Q24").Select()
